$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.001999999999999995
$ws.Range("E2").Value = 0.425
$ws.Range("G2").Value = -0.1087695748390715
$ws.Range("H2").Value = -0.1087695748390715
$ws.Range("I2").Value = -0.1002685971321625
$ws.Range("J2").Value = -0.09853982821609077
$ws.Range("K2").Value = -160.589
$ws.Range("L2").Value = -0.1884312202402604
$ws.Range("M2").Value = 0.051
$ws.Range("N2").Value = 0.0005528635077563499
$ws.Range("O2").Value = -0.0003175809052923924
$ws.Range("P2").Value = 0.051
$ws.Range("Q2").Value = 0.0005528635077563499
$ws.Range("R2").Value = -0.0003175809052923924
$ws.Range("U2").Value = 31.376
$ws.Range("V2").Value = 0.3401303023404555
$ws.Range("W2").Value = -0.8193614859546943
$ws.Range("X2").Value = 0.1827459789052231
$ws.Range("Y2").Value = -1.002107464859917
$ws.Range("Z2").Value = 3.685967757703689
$ws.Range("AA2").Value = -0.143843224180298
$ws.Range("AB2").Value = 0.1055253777724187
$ws.Range("AC2").Value = -0.2493686019527167
$ws.Range("AD2").Value = 241.8
$ws.Range("AE2").Value = 17.76554878554231
$ws.Range("AF2").Value = 259.5655487855423
$ws.Range("AG2").Value = 228.1895487855423
$ws.Range("AH2").Value = 0.7377950265889126
$ws.Range("AI2").Value = 1.454918586066854
$ws.Range("AJ2").Value = 0.7121208540361047
$ws.Range("AK2").Value = 1.55199788525761
$ws.Range("AL2").Value = 50.3
$ws.Range("AM2").Value = 47.51
$ws.Range("AN2").Value = -3.312464895817636
$ws.Range("AO2").Value = -1.68986083499006
$ws.Range("AP2").Value = -3.126012696214122
$ws.Range("AQ2").Value = -1.789097032203747

# Row 3
$ws.Range("D3").Value = -0.103
$ws.Range("E3").Value = 0.425
$ws.Range("G3").Value = 0.7183098591549296
$ws.Range("H3").Value = 0.7183098591549296
$ws.Range("I3").Value = 0.7042253521126761
$ws.Range("J3").Value = 0.6799417192812045
$ws.Range("K3").Value = 0.111
$ws.Range("L3").Value = 0.7816901408450705
$ws.Range("M3").Value = 0.017
$ws.Range("N3").Value = 0.03107861060329068
$ws.Range("O3").Value = 0.1531531531531532
$ws.Range("P3").Value = 0.017
$ws.Range("Q3").Value = 0.03107861060329068
$ws.Range("R3").Value = 0.1531531531531532
$ws.Range("U3").Value = 0.476
$ws.Range("V3").Value = 0.8702010968921389
$ws.Range("W3").Value = 0.07449664429530202
$ws.Range("X3").Value = 0.09963911376479914
$ws.Range("Y3").Value = -0.02514246946949712
$ws.Range("Z3").Value = 0.1238012205754141
$ws.Range("AA3").Value = 0.08417761476715871
$ws.Range("AB3").Value = 0.09963911376479914
$ws.Range("AC3").Value = -0.01546149899764043
$ws.Range("AG3").Value = -0.476
$ws.Range("AJ3").Value = -6.70422535211267
$ws.Range("AK3").Value = -0.4937759336099585
$ws.Range("AP3").Value = -4.621359223300971

# Row 4
$ws.Range("D4").Value = 0.107
$ws.Range("G4").Value = -0.1089074052341274
$ws.Range("H4").Value = -0.1089074052341274
$ws.Range("I4").Value = -0.1004026637215215
$ws.Range("J4").Value = -0.1004026637215215
$ws.Range("K4").Value = -160.7
$ws.Range("L4").Value = -0.1885928881586668
$ws.Range("M4").Value = 0.034
$ws.Range("N4").Value = 0.0003707742639040349
$ws.Range("O4").Value = -0.0002115743621655259
$ws.Range("P4").Value = 0.034
$ws.Range("Q4").Value = 0.0003707742639040349
$ws.Range("R4").Value = -0.0002115743621655259
$ws.Range("U4").Value = 30.9
$ws.Range("V4").Value = 0.3369683751363141
$ws.Range("W4").Value = -1.713219616204691
$ws.Range("X4").Value = 0.2658528440456471
$ws.Range("Y4").Value = -1.979072460250338
$ws.Range("Z4").Value = 3.703727066038439
$ws.Range("AA4").Value = -0.3718640631277548
$ws.Range("AB4").Value = 0.1114116417800382
$ws.Range("AC4").Value = -0.483275704907793
$ws.Range("AD4").Value = 241.8
$ws.Range("AE4").Value = 17.76554878554231
$ws.Range("AF4").Value = 259.5655487855423
$ws.Range("AG4").Value = 228.6655487855423
$ws.Range("AH4").Value = 0.7389439405115545
$ws.Range("AI4").Value = 1.466757516176777
$ws.Range("AJ4").Value = 0.7137644782729574
$ws.Range("AK4").Value = 1.565499535563144
$ws.Range("AL4").Value = 50.3
$ws.Range("AM4").Value = 47.51
$ws.Range("AN4").Value = -3.307797537619699
$ws.Range("AO4").Value = -1.691848906560636
$ws.Range("AP4").Value = -3.128119682428759
$ws.Range("AQ4").Value = -1.791201852241633

